$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2023-07-18) was added to the data set. It is
# inserted right after the existing row 14, which pushes every following
# row (old rows 15-48) down by one (new rows 16-49).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's values.
$ws.Cells.Item(15,1).Value  = 1
$ws.Cells.Item(15,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(15,4).Value  = 45125
$ws.Cells.Item(15,5).Value  = 15
$ws.Cells.Item(15,6).Value  = 100112044
$ws.Cells.Item(15,7).Value  = "Perejil"
$ws.Cells.Item(15,8).Value  = "Sin especificar"
$ws.Cells.Item(15,9).Value  = "Primera"
$ws.Cells.Item(15,10).Value = 350
$ws.Cells.Item(15,11).Value = 800
$ws.Cells.Item(15,12).Value = 1000
$ws.Cells.Item(15,13).Value = 857
$ws.Cells.Item(15,14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(15,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15,16).Value = 428
$ws.Cells.Item(15,17).Value = 2
$ws.Cells.Item(15,18).Value = "Hortaliza"
